$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update PO numbers in column C for rows 2-4
$ws.Range("C2").Value = 1774017
$ws.Range("C3").Value = 1774018
$ws.Range("C4").Value = 1774019

# Clear the values in row 5 (B5 and C5) but keep formatting/styles
$ws.Range("B5:C5").ClearContents()

# Update the active selection to B5
$ws.Activate()
$ws.Range("B5").Select()
